# HU1-CP001 "Validating ux": rename the "visibility" validation keyword to
# "visible" across the test-data sheets, and switch the active/selected
# worksheet from "body" back to "menu" (with a full-column selection on
# column D).

$wb = $excel.ActiveWorkbook

$wsMenu = $wb.Worksheets.Item("menu")
$wsBody = $wb.Worksheets.Item("body")

# --- Column D ("Validate") values: "visibility" -> "visible" -------------
# menu sheet: data rows 2..7
for ($r = 2; $r -le 7; $r++) {
    $cell = $wsMenu.Cells.Item($r, 4)
    if ($cell.Value() -eq "visibility") {
        $cell.Value = "visible"
    }
}

# body sheet: data rows 2..8 (row 6 was previously blank and now gets the
# value too)
for ($r = 2; $r -le 8; $r++) {
    $wsBody.Cells.Item($r, 4).Value = "visible"
}

# --- Sheet selection / active tab -----------------------------------------
# Previously "body" was the active/selected tab with D2 selected.
# Now "menu" becomes active/selected, with D2:D7 selected, and "body"
# reverts to its non-selected state (its own selection stays D2).
$wsMenu.Activate()
$wsMenu.Range("D2:D7").Select()
